$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: B2/C2 become numeric 0 (were the shared string "-") ---
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0

# --- New column G: per-segment Area, and H2: Atotal (sum) ---
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Summary cells J2/K2 (with headers Atotal / Qtotal) ---
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- Selection matches the saved view state ---
$null = $ws.Range("J2:K2").Select()
